# Applies the cryptos-list price/volume refresh described by the commit
# "Updated cryptos list on Thu Sep 21 05:22:19 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target cell + new text. "Force" marks values that look like
# plain numbers (e.g. "214.41") but must stay text, matching the original
# inline-string cells -- same trick Excel's UI uses (leading apostrophe).
$updates = @(
    @{ Cell = "D2"; Value = '27.082.69'; Force = $False },
    @{ Cell = "E2"; Value = '  -0.09%  '; Force = $False },
    @{ Cell = "D3"; Value = '1.623.99'; Force = $False },
    @{ Cell = "E3"; Value = '  -0.86%  '; Force = $False },
    @{ Cell = "E4"; Value = '  -0.20%  '; Force = $False },
    @{ Cell = "D5"; Value = '214.41'; Force = $True },
    @{ Cell = "E5"; Value = '  -1.01%  '; Force = $False },
    @{ Cell = "D6"; Value = '0.517'; Force = $True },
    @{ Cell = "E6"; Value = '  -0.79%  '; Force = $False },
    @{ Cell = "E7"; Value = '  -0.21%  '; Force = $False },
    @{ Cell = "D8"; Value = '0.0630'; Force = $True },
    @{ Cell = "E8"; Value = '  +0.83%  '; Force = $False },
    @{ Cell = "D9"; Value = '0.250'; Force = $True },
    @{ Cell = "E9"; Value = '  -1.49%  '; Force = $False },
    @{ Cell = "D10"; Value = '20.02'; Force = $True },
    @{ Cell = "E10"; Value = '  +0.55%  '; Force = $False },
    @{ Cell = "D11"; Value = '0.0847'; Force = $True },
    @{ Cell = "E11"; Value = '  +0.02%  '; Force = $False },
    @{ Cell = "D12"; Value = '1.852.26'; Force = $False },
    @{ Cell = "E12"; Value = '  -0.89%  '; Force = $False },
    @{ Cell = "D13"; Value = '1.643.45'; Force = $False },
    @{ Cell = "E13"; Value = '  +0.27%  '; Force = $False },
    @{ Cell = "D14"; Value = '4.14'; Force = $True },
    @{ Cell = "E14"; Value = '  +0.27%  '; Force = $False },
    @{ Cell = "D15"; Value = '0.540'; Force = $True },
    @{ Cell = "E15"; Value = '  -0.21%  '; Force = $False },
    @{ Cell = "D16"; Value = '64.59'; Force = $True },
    @{ Cell = "E16"; Value = '  -3.26%  '; Force = $False },
    @{ Cell = "D17"; Value = '27.065.67'; Force = $False },
    @{ Cell = "E17"; Value = '  -0.27%  '; Force = $False },
    @{ Cell = "D18"; Value = '0.0₃0742'; Force = $False },
    @{ Cell = "E18"; Value = '  +0.31%  '; Force = $False },
    @{ Cell = "D19"; Value = '213.41'; Force = $True },
    @{ Cell = "E19"; Value = '  -1.53%  '; Force = $False },
    @{ Cell = "E20"; Value = '  -0.14%  '; Force = $False },
    @{ Cell = "D21"; Value = '6.82'; Force = $True },
    @{ Cell = "E21"; Value = '  -1.74%  '; Force = $False },
    @{ Cell = "E22"; Value = '  -1.29%  '; Force = $False },
    @{ Cell = "D23"; Value = '2.34'; Force = $True },
    @{ Cell = "E23"; Value = '  -7.63%  '; Force = $False },
    @{ Cell = "D24"; Value = '9.07'; Force = $True },
    @{ Cell = "E24"; Value = '  -0.42%  '; Force = $False },
    @{ Cell = "D25"; Value = '148.03'; Force = $True },
    @{ Cell = "E25"; Value = '  +0.85%  '; Force = $False },
    @{ Cell = "E26"; Value = '  -0.22%  '; Force = $False },
    @{ Cell = "E27"; Value = '  -0.50%  '; Force = $False },
    @{ Cell = "E28"; Value = '  -3.04%  '; Force = $False },
    @{ Cell = "D29"; Value = '15.59'; Force = $True },
    @{ Cell = "E29"; Value = '  -0.45%  '; Force = $False },
    @{ Cell = "E30"; Value = '  +0.76%  '; Force = $False },
    @{ Cell = "E31"; Value = '  -0.87%  '; Force = $False },
    @{ Cell = "E32"; Value = '  -0.96%  '; Force = $False },
    @{ Cell = "D33"; Value = '0.732'; Force = $True },
    @{ Cell = "E33"; Value = '  +34.50%  '; Force = $False },
    @{ Cell = "D34"; Value = '3.00'; Force = $True },
    @{ Cell = "E34"; Value = '  -0.49%  '; Force = $False },
    @{ Cell = "D35"; Value = '1.362.22'; Force = $False },
    @{ Cell = "E35"; Value = '  +4.24%  '; Force = $False },
    @{ Cell = "E36"; Value = '  +0.48%  '; Force = $False },
    @{ Cell = "E37"; Value = '  -0.84%  '; Force = $False },
    @{ Cell = "E38"; Value = '  +0.85%  '; Force = $False },
    @{ Cell = "D39"; Value = '0.843'; Force = $True },
    @{ Cell = "E39"; Value = '  -1.75%  '; Force = $False },
    @{ Cell = "B41"; Value = 'MXToken'; Force = $False },
    @{ Cell = "C41"; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; Force = $False },
    @{ Cell = "D41"; Value = '2.23'; Force = $True },
    @{ Cell = "E41"; Value = '  +0.21%  '; Force = $False },
    @{ Cell = "B42"; Value = 'TrustWalletToken'; Force = $False },
    @{ Cell = "C42"; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; Force = $False },
    @{ Cell = "D42"; Value = '0.799'; Force = $True },
    @{ Cell = "E42"; Value = '  -1.36%  '; Force = $False },
    @{ Cell = "D43"; Value = '64.59'; Force = $True },
    @{ Cell = "E43"; Value = '  +4.74%  '; Force = $False },
    @{ Cell = "D44"; Value = '5.33'; Force = $True },
    @{ Cell = "E44"; Value = '  +0.75%  '; Force = $False },
    @{ Cell = "D45"; Value = '1.763.26'; Force = $False },
    @{ Cell = "E45"; Value = '  -0.91%  '; Force = $False },
    @{ Cell = "D46"; Value = '1.65'; Force = $True },
    @{ Cell = "E46"; Value = '  +3.22%  '; Force = $False },
    @{ Cell = "D47"; Value = '89.87'; Force = $True },
    @{ Cell = "E47"; Value = '  -1.73%  '; Force = $False },
    @{ Cell = "D48"; Value = '0.863'; Force = $True },
    @{ Cell = "E48"; Value = '  +29.26%  '; Force = $False },
    @{ Cell = "E49"; Value = '  -1.04%  '; Force = $False },
    @{ Cell = "E50"; Value = '  +5.05%  '; Force = $False },
    @{ Cell = "D51"; Value = '0.0513'; Force = $True },
    @{ Cell = "E51"; Value = '  +0.36%  '; Force = $False }
)

foreach ($u in $updates) {
    $value = $u.Value
    if ($u.Force) {
        # leading apostrophe forces text storage instead of numeric auto-conversion
        $value = "'" + $value
    }
    $ws.Range($u.Cell).Value = $value
}
